# Rename the "results" item header and its corresponding jxls expression
# to "Item" / ${results.ITEM_NAME}, reflecting the new sample jxls query
# (with chart) in the demo database, and move the sheet selection to C11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = "Item"
$ws.Range("C10").Value = '${results.ITEM_NAME}'

$ws.Range("C11").Select()
